# Auto-generated script applying Typhon_Profits market-value refresh
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1058.1364
$ws.Range("I40").Value = 791.6667
$ws.Range("J40").Value = 1377.9
$ws.Range("K40").Value = 791.6667
$ws.Range("L40").Value = 1377.9
$ws.Range("M40").Value = -616.6667
$ws.Range("N40").Value = -1727.9
$ws.Range("H74").Value = 7816999
$ws.Range("I74").Value = 4283.5713
$ws.Range("J74").Value = 13893556
$ws.Range("K74").Value = 4283.5713
$ws.Range("L74").Value = 13893556
$ws.Range("M74").Value = -3347.5713
$ws.Range("N74").Value = -13895428
$ws.Range("H77").Value = 7816999
$ws.Range("I77").Value = 4283.5713
$ws.Range("J77").Value = 13893556
$ws.Range("K77").Value = 21417.8565
$ws.Range("L77").Value = 69467780
$ws.Range("M77").Value = -16737.8565
$ws.Range("N77").Value = -69477140
$ws.Range("H100").Value = 2222.6667
$ws.Range("I100").Value = 1560.8
$ws.Range("K100").Value = 1560.8
$ws.Range("M100").Value = -1019.8
$ws.Range("H129").Value = 759.5454999999999
$ws.Range("J129").Value = 821.7368
$ws.Range("L129").Value = 2465.2104
$ws.Range("N129").Value = -12465.2104
$ws.Range("H137").Value = 104338.46
$ws.Range("I137").Value = 122433.24
$ws.Range("K137").Value = 367299.72
$ws.Range("M137").Value = -364749.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2899.4167
$ws.Range("I45").Value = 3083.25
$ws.Range("J45").Value = 2715.5833
$ws.Range("K45").Value = 3083.25
$ws.Range("L45").Value = 2715.5833
$ws.Range("M45").Value = -2706.25
$ws.Range("N45").Value = -3469.5833
$ws.Range("H61").Value = 3353.1143
$ws.Range("I61").Value = 3100.4138
$ws.Range("K61").Value = 3100.4138
$ws.Range("M61").Value = -2888.4138
$ws.Range("H136").Value = 3353.1143
$ws.Range("I136").Value = 3100.4138
$ws.Range("K136").Value = 9301.241399999999
$ws.Range("M136").Value = -6751.241399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2019.9166
$ws.Range("I20").Value = 1870.8182
$ws.Range("K20").Value = 1870.8182
$ws.Range("M20").Value = -1623.8182
$ws.Range("H80").Value = 732.28125
$ws.Range("I80").Value = 774.5
$ws.Range("J80").Value = 706.95
$ws.Range("K80").Value = 774.5
$ws.Range("L80").Value = 706.95
$ws.Range("M80").Value = 223.5
$ws.Range("N80").Value = -2702.95
$ws.Range("H83").Value = 732.28125
$ws.Range("I83").Value = 774.5
$ws.Range("J83").Value = 706.95
$ws.Range("K83").Value = 3872.5
$ws.Range("L83").Value = 3534.75
$ws.Range("M83").Value = 1119.5
$ws.Range("N83").Value = -13518.75
$ws.Range("H86").Value = 1773.5862
$ws.Range("I86").Value = 1638.0555
$ws.Range("J86").Value = 1995.3636
$ws.Range("K86").Value = 1638.0555
$ws.Range("L86").Value = 1995.3636
$ws.Range("M86").Value = -515.0554999999999
$ws.Range("N86").Value = -4241.3636
$ws.Range("H89").Value = 1773.5862
$ws.Range("I89").Value = 1638.0555
$ws.Range("J89").Value = 1995.3636
$ws.Range("K89").Value = 8190.2775
$ws.Range("L89").Value = 9976.817999999999
$ws.Range("M89").Value = -2574.2775
$ws.Range("N89").Value = -21208.818
$ws.Range("H105").Value = 1001966
$ws.Range("I105").Value = 1600.4762
$ws.Range("K105").Value = 1600.4762
$ws.Range("M105").Value = 146.5237999999999
$ws.Range("H126").Value = 58500
$ws.Range("J126").Value = 58500
$ws.Range("L126").Value = 58500
$ws.Range("N126").Value = -68380

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 141.08333
$ws.Range("I22").Value = 139.73334
$ws.Range("J22").Value = 143.33333
$ws.Range("K22").Value = 139.73334
$ws.Range("L22").Value = 143.33333
$ws.Range("M22").Value = 210.26666
$ws.Range("N22").Value = -843.3333299999999
$ws.Range("H31").Value = 4795
$ws.Range("I31").Value = 2104.35
$ws.Range("J31").Value = 7960.4707
$ws.Range("K31").Value = 2104.35
$ws.Range("L31").Value = 7960.4707
$ws.Range("M31").Value = -1809.35
$ws.Range("N31").Value = -8550.4707
$ws.Range("H34").Value = 4795
$ws.Range("I34").Value = 2104.35
$ws.Range("J34").Value = 7960.4707
$ws.Range("K34").Value = 2104.35
$ws.Range("L34").Value = 7960.4707
$ws.Range("M34").Value = -1902.35
$ws.Range("N34").Value = -8364.4707
$ws.Range("H141").Value = 28330.145
$ws.Range("J141").Value = 28330.145
$ws.Range("L141").Value = 28330.145
$ws.Range("N141").Value = -38690.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 10088.5
$ws.Range("I2").Value = 16697.5
$ws.Range("J2").Value = 175
$ws.Range("K2").Value = 100185
$ws.Range("L2").Value = 1050
$ws.Range("M2").Value = -100072
$ws.Range("N2").Value = -1276
$ws.Range("H19").Value = 1549.5
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 9000
$ws.Range("N19").Value = -9348
$ws.Range("H131").Value = 726.67
$ws.Range("J131").Value = 794.34485
$ws.Range("L131").Value = 2383.03455
$ws.Range("N131").Value = -12463.03455
$ws.Range("H137").Value = 13893153
$ws.Range("I137").Value = 1006.6667
$ws.Range("J137").Value = 22228440
$ws.Range("K137").Value = 3020.0001
$ws.Range("L137").Value = 66685320
$ws.Range("M137").Value = 2079.9999
$ws.Range("N137").Value = -66695520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 69.52941
$ws.Range("I2").Value = 73.38461
$ws.Range("K2").Value = 73.38461
$ws.Range("M2").Value = 39.61539

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2732.6667
$ws.Range("I100").Value = 2399
$ws.Range("J100").Value = 2799.4
$ws.Range("K100").Value = 2399
$ws.Range("L100").Value = 2799.4
$ws.Range("M100").Value = -1858
$ws.Range("N100").Value = -3881.4
$ws.Range("H122").Value = 703422.4
$ws.Range("I122").Value = 1035201.44
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 3105604.32
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3103154.32
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4333.3335
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4333.3335
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -31240
$ws.Range("H81").Value = 4685
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 4685
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H122").Value = 1256.25
$ws.Range("I122").Value = 1007.1429
$ws.Range("K122").Value = 3021.4287
$ws.Range("M122").Value = -571.4287000000004
